$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "28.422.77"
$ws.Range("E2").Value = "  +0.11%  "

# Row 3
$ws.Range("D3").Value = "1.816.83"
$ws.Range("E3").Value = "  -0.41%  "

# Row 4
$ws.Range("E4").Value = "  +0.17%  "

# Row 5
$ws.Range("D5").Value = "'315.03"
$ws.Range("E5").Value = "  -0.77%  "

# Row 6
$ws.Range("E6").Value = "  +0.11%  "

# Row 7
$ws.Range("D7").Value = "'0.5093"
$ws.Range("E7").Value = "  -4.40%  "

# Row 8
$ws.Range("D8").Value = "'0.3949"
$ws.Range("E8").Value = "  -2.23%  "

# Row 9
$ws.Range("D9").Value = "'0.08058"
$ws.Range("E9").Value = "  +5.72%  "

# Row 10
$ws.Range("D10").Value = "'41.68"
$ws.Range("E10").Value = "  -0.28%  "

# Row 11
$ws.Range("D11").Value = "'1.106"
$ws.Range("E11").Value = "  -0.25%  "

# Row 12
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").Value = "'6.283"
$ws.Range("E12").Value = "  -0.65%  "

# Row 13
$ws.Range("B13").Value = "Solana"
$ws.Range("C13").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D13").Value = "'20.97"
$ws.Range("E13").Value = "  +0.54%  "

# Row 14
$ws.Range("E14").Value = "  +0.13%  "

# Row 15
$ws.Range("D15").Value = "'7.494"
$ws.Range("E15").Value = "  -1.39%  "

# Row 16
$ws.Range("D16").Value = "1.819.50"
$ws.Range("E16").Value = "  -0.50%  "

# Row 17
$ws.Range("D17").Value = "'0.00001132"
$ws.Range("E17").Value = "  +5.19%  "

# Row 18
$ws.Range("D18").Value = "'92.52"
$ws.Range("E18").Value = "  +3.44%  "

# Row 19
$ws.Range("D19").Value = "'0.06642"
$ws.Range("E19").Value = "  +0.68%  "

# Row 20
$ws.Range("D20").Value = "'17.68"
$ws.Range("E20").Value = "  +0.03%  "

# Row 21
$ws.Range("E21").Value = "  +0.10%  "

# Row 22
$ws.Range("D22").Value = "'6.087"
$ws.Range("E22").Value = "  +0.13%  "

# Row 23
$ws.Range("D23").Value = "28.451.02"
$ws.Range("E23").Value = "  +0.19%  "

# Row 24
$ws.Range("E24").Value = "  +0.87%  "

# Row 25
$ws.Range("D25").Value = "'2.273"
$ws.Range("E25").Value = "  +3.09%  "

# Row 26
$ws.Range("D26").Value = "'21.13"
$ws.Range("E26").Value = "  +2.48%  "

# Row 27
$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").Value = "'155.84"
$ws.Range("E27").Value = "  -1.10%  "

# Row 28
$ws.Range("B28").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C28").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D28").Value = "2.027.45"
$ws.Range("E28").Value = "  -0.56%  "

# Row 29
$ws.Range("D29").Value = "'2.401"
$ws.Range("E29").Value = "  -2.11%  "

# Row 30
$ws.Range("D30").Value = "'125.93"
$ws.Range("E30").Value = "  +1.57%  "

# Row 31
$ws.Range("D31").Value = "'0.1099"
$ws.Range("E31").Value = "  -0.30%  "

# Row 32
$ws.Range("D32").Value = "'1.106"
$ws.Range("E32").Value = "  -1.48%  "

# Row 33
$ws.Range("D33").Value = "'5.777"
$ws.Range("E33").Value = "  +2.13%  "

# Row 34
$ws.Range("E34").Value = "  +0.16%  "

# Row 35
$ws.Range("D35").Value = "'0.07010"
$ws.Range("E35").Value = "  -5.44%  "

# Row 36
$ws.Range("D36").Value = "'0.2220"
$ws.Range("E36").Value = "  -0.57%  "

# Row 37
$ws.Range("D37").Value = "'5.218"
$ws.Range("E37").Value = "  +0.39%  "

# Row 38
$ws.Range("D38").Value = "'0.02326"
$ws.Range("E38").Value = "  -0.68%  "

# Row 39
$ws.Range("D39").Value = "'8.815"
$ws.Range("E39").Value = "  -0.97%  "

# Row 40
$ws.Range("D40").Value = "'0.6283"
$ws.Range("E40").Value = "  +0.49%  "

# Row 41
$ws.Range("D41").Value = "'11.28"
$ws.Range("E41").Value = "  -0.23%  "

# Row 42
$ws.Range("E42").Value = "  -0.20%  "

# Row 43
$ws.Range("D43").Value = "'1.002"
$ws.Range("E43").Value = "  +0.14%  "

# Row 44
$ws.Range("D44").Value = "'1.400"
$ws.Range("E44").Value = "  +0.40%  "

# Row 45
$ws.Range("D45").Value = "'13.47"
$ws.Range("E45").Value = "  -0.50%  "

# Row 46
$ws.Range("D46").Value = "'3.740"
$ws.Range("E46").Value = "  +1.07%  "

# Row 47
$ws.Range("D47").Value = "'0.5912"
$ws.Range("E47").Value = "  +1.29%  "

# Row 48
$ws.Range("D48").Value = "'124.84"
$ws.Range("E48").Value = "  -0.11%  "

# Row 49
$ws.Range("D49").Value = "'1.971"
$ws.Range("E49").Value = "  -0.97%  "

# Row 50
$ws.Range("D50").Value = "'1.186"
$ws.Range("E50").Value = "  -1.10%  "

# Row 51
$ws.Range("D51").Value = "'0.06885"
$ws.Range("E51").Value = "  -0.08%  "
